$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "Data" to "Summary"
$ws.Name = "Summary"

# Re-assert formatting on the untouched header cells (A1 "name" style,
# A3 "title" style) so their appearance is preserved through the edit.
$ws.Range("A1").Font.Size = 18
$ws.Range("A3").Font.Bold = $true

# Remove the old block of cells (rows 5-8); the content is relocated further
# down the sheet and a new "Source Type" note is inserted above it.
$ws.Range("A5:D8").Clear()

# New "Source Type" note (bold + underlined) at row 9
$ws.Range("A9").Value = "Source Type: SME Associations (Most Widely Used)"
$ws.Range("A9").Font.Bold = $true
$ws.Range("A9").Font.Underline = $true

# Column headers (bold "title" style) now at row 11
$ws.Range("B11").Value = "Micro"
$ws.Range("C11").Value = "SMEs"
$ws.Range("D11").Value = "MSMEs"
$ws.Range("B11:D11").Font.Bold = $true

# Data row "Enterprises (absolute #)" now at row 12
$ws.Range("A12").Value = "Enterprises (absolute #)"
$ws.Range("A12").Font.Bold = $true

$numCells12 = @("B12", "C12", "D12")
$numVals12 = @("19371", "783", "20154")
for ($i = 0; $i -lt $numCells12.Length; $i++) {
    $c = $ws.Range($numCells12[$i])
    $c.NumberFormat = "@"
    $c.Value = $numVals12[$i]
    $c.Style = "Normal"
}

# Data row "Enterprises density (per 1000 people)" now at row 13
$ws.Range("A13").Value = "Enterprises density (per 1000 people)"
$ws.Range("A13").Font.Bold = $true

$numCells13 = @("B13", "C13", "D13")
$numVals13 = @("1", "0", "1")
for ($i = 0; $i -lt $numCells13.Length; $i++) {
    $c = $ws.Range($numCells13[$i])
    $c.NumberFormat = "@"
    $c.Value = $numVals13[$i]
    $c.Style = "Normal"
}

# Source note (italic "source" style) now at row 14
$ws.Range("A14").Value = "Source: MDE, 2011"
$ws.Range("A14").Font.Italic = $true

# New "MDE" title (bold "title" style) at row 20
$ws.Range("A20").Value = "MDE"
$ws.Range("A20").Font.Bold = $true

# New citation (italic "source" style) at row 21
$ws.Range("A21").Value = 'Ministeriio da Economia (MDE), "S' + [char]0x00ED + 'ntese do Programa de Desenvolvimento das MPME`s", 2012, p. 6. Avaialable at http://www.minec.gov.ao/VerPublicacao.aspx?id=820'
$ws.Range("A21").Font.Italic = $true
